# Update country/province COVID-19 stats and re-sort by total cases
# (descending), matching the "Update countries & provincias Spain" refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp of this data refresh (cell A1).
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 22:20"

# Full refreshed table: Pais, Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes - already sorted
# descending by "Casos totales" (column B), rows 4-210.
$data = @(
    @("Estados Unidos", 240344, 25341, 10365, 224172, 5421, 705, 5807),
    @("Italia", 115242, 4668, 18278, 83049, 4053, 760, 13915),
    @("España", 110238, 6120, 26743, 73399, 6092, 709, 10096),
    @("Alemania", 84688, 6707, 21400, 62184, 3936, 173, 1104),
    @("China", 81589, 35, 76408, 1863, 429, 6, 3318),
    @("Francia", 59105, 2116, 12428, 42174, 6399, 471, 4503),
    @("Iran", 50468, 2875, 16711, 30597, 3956, 124, 3160),
    @("Reino Unido", 33718, 4244, 135, 30662, 163, 569, 2921),
    @("Suiza", 18827, 1059, 4013, 14278, 348, 48, 536),
    @("Turquia", 18135, 2456, 415, 17364, 1101, 79, 356),
    @("Belgica", 15348, 1384, 2495, 11842, 1144, 183, 1011),
    @("Paises Bajos", 14697, 1083, 250, 13108, 1053, 166, 1339),
    @("Canada", 11131, 1400, 1906, 9064, 120, 47, 161),
    @("Austria", 11108, 397, 1749, 9201, 227, 12, 158),
    @("Corea del Sur", 9976, 89, 5828, 3979, 55, 4, 169),
    @("Portugal", 9034, 783, 68, 8757, 230, 22, 209),
    @("Brasil", 7910, 1030, 127, 7484, 296, 57, 299),
    @("Israel", 6857, 765, 338, 6483, 107, 10, 36),
    @("Suecia", 5568, 621, 103, 5157, 429, 69, 308),
    @("Noruega", 5141, 264, 32, 5059, 96, 6, 50),
    @("Australia", 5139, 91, 345, 4767, 50, 4, 27),
    @("Irlanda", 3849, 402, 5, 3746, 109, 13, 98),
    @("Chequia", 3805, 216, 67, 3694, 72, 5, 44),
    @("Rusia", 3548, 771, 235, 3283, 8, 6, 30),
    @("Chile", 3404, 373, 335, 3051, 31, 2, 18),
    @("Dinamarca", 3386, 279, 1089, 2174, 153, 19, 123),
    @("Ecuador", 3163, 405, 58, 2985, 100, 22, 120),
    @("Malasia", 3116, 208, 767, 2299, 105, 5, 50),
    @("Polonia", 2946, 392, 56, 2833, 50, 14, 57),
    @("Rumania", 2738, 278, 267, 2356, 78, 23, 115),
    @("Filipinas", 2633, 322, 51, 2475, 1, 11, 107),
    @("India", 2536, 538, 189, 2275, 0, 14, 72),
    @("Luxemburgo", 2487, 168, 80, 2377, 31, 1, 30),
    @("Pakistan", 2386, 268, 107, 2245, 9, 7, 34),
    @("Japon", 2384, 0, 472, 1855, 69, 0, 57),
    @("Arabia Saudita", 1885, 165, 328, 1536, 31, 5, 21),
    @("Tailandia", 1875, 104, 505, 1355, 23, 3, 15),
    @("Indonesia", 1790, 113, 112, 1508, 0, 13, 170),
    @("Grecia", 1544, 129, 61, 1430, 91, 2, 53),
    @("Finlandia", 1518, 72, 300, 1199, 65, 2, 19),
    @("Sudafrica", 1462, 82, 50, 1407, 7, 0, 5),
    @("Peru", 1414, 91, 537, 822, 51, 17, 55),
    @("Republica Dominicana", 1380, 96, 16, 1304, 147, 3, 60),
    @("Mexico", 1378, 163, 35, 1306, 1, 8, 37),
    @("Islandia", 1319, 99, 284, 1031, 12, 2, 4),
    @("Panama", 1317, 0, 9, 1276, 50, 0, 32),
    @("Serbia", 1171, 111, 42, 1098, 81, 3, 31),
    @("Colombia", 1161, 96, 55, 1087, 50, 2, 19),
    @("Argentina", 1133, 0, 256, 843, 0, 2, 34),
    @("Singapur", 1049, 49, 266, 779, 24, 1, 4),
    @("Emiratos Arabes Unidos", 1024, 210, 96, 920, 2, 0, 8),
    @("Croacia", 1011, 48, 88, 916, 34, 1, 7),
    @("Argelia", 986, 139, 61, 839, 0, 28, 86),
    @("Catar", 949, 114, 72, 874, 37, 1, 3),
    @("Ucrania", 897, 103, 19, 856, 0, 2, 22),
    @("Eslovenia", 897, 56, 70, 810, 31, 2, 17),
    @("Egipto", 865, 86, 201, 606, 0, 6, 58),
    @("Estonia", 858, 79, 45, 802, 16, 6, 11),
    @("Hong Kong", 802, 36, 154, 644, 8, 0, 4),
    @("Nueva Zelanda", 797, 89, 92, 704, 2, 0, 1),
    @("Irak", 772, 44, 202, 516, 0, 2, 54),
    @("Crucero", 712, 0, 619, 82, 15, 0, 11),
    @("Marruecos", 691, 37, 30, 617, 1, 5, 44),
    @("Armenia", 663, 92, 33, 623, 30, 3, 7),
    @("Lituania", 649, 68, 7, 633, 11, 1, 9),
    @("Barein", 643, 74, 381, 258, 3, 0, 4),
    @("Hungria", 585, 60, 42, 522, 17, 1, 21),
    @("Bosnia y Herzegovina", 532, 73, 20, 496, 4, 3, 16),
    @("Moldavia", 505, 82, 23, 476, 65, 1, 6),
    @("Libano", 494, 15, 46, 432, 3, 2, 16),
    @("Letonia", 458, 12, 31, 427, 3, 0, 0),
    @("Bulgaria", 457, 35, 25, 422, 18, 0, 10),
    @("Tunez", 455, 32, 5, 436, 10, 2, 14),
    @("Kazajistan", 435, 55, 27, 405, 6, 0, 3),
    @("Principado de Andorra", 428, 38, 10, 403, 12, 1, 15),
    @("Eslovaquia", 426, 26, 5, 420, 3, 0, 1),
    @("Azerbaiyan", 400, 41, 26, 369, 7, 0, 5),
    @("Costa Rica", 396, 21, 6, 388, 11, 0, 2),
    @("Republica de Macedonia", 384, 30, 17, 356, 4, 0, 11),
    @("Republica de Chipre", 356, 36, 28, 318, 11, 1, 10),
    @("Uruguay", 350, 0, 62, 284, 15, 2, 4),
    @("Kuwait", 342, 25, 81, 261, 15, 0, 0),
    @("Taiwan", 339, 10, 50, 284, 0, 0, 5),
    @("Reunion", 308, 27, 40, 268, 3, 0, 0),
    @("Camerun", 306, 73, 10, 289, 0, 1, 7),
    @("Bielorrusia", 304, 141, 53, 247, 11, 2, 4),
    @("Jordania", 299, 21, 45, 249, 5, 0, 5),
    @("Burkina Faso", 288, 6, 50, 222, 0, 0, 16),
    @("Albania", 277, 18, 67, 194, 7, 1, 16),
    @("Afganistan", 273, 36, 10, 257, 0, 2, 6),
    @("San Marino", 245, 9, 21, 194, 15, 2, 30),
    @("Cuba", 233, 21, 13, 214, 7, 0, 6),
    @("Oman", 231, 21, 57, 173, 3, 0, 1),
    @("Vietnam", 227, 9, 75, 152, 3, 0, 0),
    @("Honduras", 219, 47, 3, 202, 4, 4, 14),
    @("Uzbekistan", 205, 24, 25, 178, 8, 0, 2),
    @("Ghana", 204, 9, 31, 168, 2, 0, 5),
    @("Malta", 196, 8, 2, 194, 2, 0, 0),
    @("Senegal", 195, 5, 55, 139, 0, 0, 1),
    @("Costa de Marfil", 194, 4, 9, 184, 0, 0, 1),
    @("Nigeria", 184, 10, 20, 162, 0, 0, 2),
    @("Islas Feroe", 177, 4, 81, 96, 1, 0, 0),
    @("Mauricio", 169, 8, 0, 162, 1, 1, 7),
    @("Estado de Palestina", 161, 27, 18, 142, 0, 0, 1),
    @("Sri Lanka", 151, 5, 21, 126, 5, 1, 4),
    @("Montenegro", 144, 21, 0, 142, 4, 0, 2),
    @("Venezuela", 144, 0, 43, 98, 6, 0, 3),
    @("Martinica", 138, 3, 27, 108, 19, 0, 3),
    @("Georgia", 134, 17, 26, 108, 6, 0, 0),
    @("Brunei", 133, 2, 56, 76, 3, 0, 1),
    @("Guadalupe", 128, 3, 24, 98, 14, 0, 6),
    @("Bolivia", 123, 8, 1, 114, 3, 1, 8),
    @("Consejo Danes para los Refugiados", 123, 14, 3, 109, 0, 2, 11),
    @("Kirguistan", 116, 5, 5, 111, 5, 0, 0),
    @("Mayotte", 116, 15, 10, 105, 3, 0, 1),
    @("Kenia", 110, 29, 4, 103, 2, 2, 3),
    @("Camboya", 110, 1, 34, 76, 1, 0, 0),
    @("Niger", 98, 24, 0, 93, 0, 0, 5),
    @("Isla de Man", 95, 27, 0, 94, 0, 0, 1),
    @("Trinidad yTobago", 94, 4, 1, 88, 0, 0, 5),
    @("Gibraltar", 88, 7, 46, 42, 0, 0, 0),
    @("Ruanda", 84, 2, 0, 84, 0, 0, 0),
    @("Paraguay", 77, 8, 2, 72, 4, 0, 3),
    @("Liechtenstein", 75, 3, 0, 75, 0, 0, 0),
    @("Aruba", 60, 5, 1, 59, 0, 0, 0),
    @("Monaco", 60, 5, 2, 57, 2, 0, 1),
    @("Madagascar", 59, 2, 0, 59, 6, 0, 0),
    @("Banglades", 56, 2, 25, 25, 1, 0, 6),
    @("Guinea", 52, 22, 0, 52, 0, 0, 0),
    @("Guayana Francesa", 51, 0, 15, 36, 0, 0, 0),
    @("Guatemala", 47, 8, 12, 34, 1, 0, 1),
    @("Barbados", 46, 1, 0, 46, 0, 0, 0),
    @("Uganda", 45, 1, 0, 45, 0, 0, 0),
    @("Jamaica", 44, 0, 2, 39, 0, 0, 3),
    @("El Salvador", 41, 8, 0, 39, 4, 0, 2),
    @("Macao", 41, 0, 10, 31, 0, 0, 0),
    @("Republica de Yibuti", 40, 7, 0, 40, 0, 0, 0),
    @("Zambia", 39, 3, 0, 38, 0, 1, 1),
    @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
    @("Togo", 39, 3, 17, 20, 0, 0, 2),
    @("Polinesia Francesa", 37, 0, 0, 37, 1, 0, 0),
    @("Mali", 36, 5, 0, 33, 0, 0, 3),
    @("Guam", 32, 0, 0, 31, 0, 0, 1),
    @("Bermudas", 32, 0, 10, 22, 0, 0, 0),
    @("Etiopia", 29, 0, 3, 26, 2, 0, 0),
    @("Eritrea", 22, 7, 0, 22, 0, 0, 0),
    @("Islas Caimanes", 22, 0, 0, 21, 0, 0, 1),
    @("San Martin (Parte Francesa)", 22, 0, 2, 19, 0, 0, 1),
    @("Congo", 22, 0, 2, 18, 0, 0, 2),
    @("Gabon", 21, 3, 0, 20, 0, 0, 1),
    @("Bahamas", 21, 0, 1, 19, 0, 0, 1),
    @("Birmania", 20, 4, 0, 19, 0, 0, 1),
    @("Tanzania", 20, 0, 2, 17, 0, 0, 1),
    @("Guyana", 19, 0, 0, 15, 0, 1, 4),
    @("Maldivas", 19, 0, 13, 6, 0, 0, 0),
    @("Nueva Caledonia", 18, 2, 1, 17, 0, 0, 0),
    @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @("Haiti", 16, 0, 1, 15, 0, 0, 0),
    @("Siria", 16, 6, 0, 14, 0, 0, 2),
    @("San Martin (Parte Holandesa)", 16, 0, 6, 9, 0, 0, 1),
    @("Guinea Ecuatorial", 15, 0, 1, 14, 0, 0, 0),
    @("Mongolia", 14, 0, 2, 12, 0, 0, 0),
    @("Namibia", 14, 0, 2, 12, 0, 0, 0),
    @("Benin", 13, 0, 1, 12, 0, 0, 0),
    @("Santa Lucia", 13, 0, 1, 12, 0, 0, 0),
    @("Dominica", 12, 0, 0, 12, 0, 0, 0),
    @("Curazao", 11, 0, 3, 7, 0, 0, 1),
    @("Mozambique", 10, 0, 0, 10, 0, 0, 0),
    @("Seychelles", 10, 0, 0, 10, 0, 0, 0),
    @("Libia", 10, 0, 0, 10, 0, 0, 0),
    @("Granada", 10, 1, 0, 10, 0, 0, 0),
    @("Laos", 10, 0, 0, 10, 0, 0, 0),
    @("Surinam", 10, 0, 0, 10, 0, 0, 0),
    @("Groenlandia", 10, 0, 2, 8, 0, 0, 0),
    @("Antigua y Barbuda", 9, 2, 0, 9, 0, 0, 0),
    @("Guinea-Bisau", 9, 0, 0, 9, 0, 0, 0),
    @("Suazilandia", 9, 0, 0, 9, 0, 0, 0),
    @("Zimbabue", 9, 1, 0, 8, 0, 0, 1),
    @("Montserrat", 9, 0, 0, 7, 0, 0, 2),
    @("San Cristobal y Nieves", 8, 0, 0, 8, 0, 0, 0),
    @("Republica del Chad", 8, 1, 0, 8, 0, 0, 0),
    @("Angola", 8, 0, 1, 5, 0, 0, 2),
    @("Sudan", 8, 1, 2, 4, 0, 0, 2),
    @("Fiyi", 7, 2, 0, 7, 0, 0, 0),
    @("Santa Sede", 6, 0, 0, 6, 0, 0, 0),
    @("Liberia", 6, 0, 0, 6, 0, 0, 0),
    @("Cabo Verde", 6, 0, 0, 5, 0, 0, 1),
    @("Nepal", 6, 1, 1, 5, 0, 0, 0),
    @("San Bartolome", 6, 0, 1, 5, 0, 0, 0),
    @("Mauritania", 6, 0, 2, 3, 0, 0, 1),
    @("Islas Turcas y Caicos", 5, 0, 0, 5, 0, 0, 0),
    @("Nicaragua", 5, 0, 0, 4, 0, 0, 1),
    @("Somalia", 5, 0, 1, 4, 0, 0, 0),
    @("Butan", 5, 1, 1, 4, 0, 0, 0),
    @("Botsuana", 4, 0, 0, 3, 0, 0, 1),
    @("Gambia", 4, 0, 2, 1, 0, 0, 1),
    @("Burundi", 3, 1, 0, 3, 0, 0, 0),
    @("Malaui", 3, 3, 0, 3, 0, 0, 0),
    @("Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0),
    @("Belice", 3, 0, 0, 3, 0, 0, 0),
    @("Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0),
    @("Anguila", 3, 1, 0, 3, 0, 0, 0),
    @("Sierra Leona", 2, 0, 0, 2, 0, 0, 0),
    @("Bonaire, San Eustaquio y Saba", 2, 0, 0, 2, 0, 0, 0),
    @("San Vicente y las Granadinas", 2, 0, 1, 1, 0, 0, 0),
    @("Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
    @("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
